# LogicComponentSequenceDiagram.pptx update
# - Refresh the "today" auto-date fields baked into the slide master /
#   slide layouts / notes master (cosmetic artifact of the deck being
#   re-saved on a later date).
# - Rename the deletePerson(p) call-out to deleteRestaurant(r) and widen
#   its text box to fit the new label.
# - Rename the saveAddressBook(AddressBook) call-out to
#   saveFoodDiary(FoodDiary).

$p = $ppt.ActivePresentation

function Update-DateShapeText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "12/15/2018") {
                $tr.Text = "4/5/19"
            }
        }
    }
}

# Notes master "today" field (best effort - some hosts treat the notes
# master as read-only, in which case this is a harmless no-op).
Update-DateShapeText $p.NotesMaster.Shapes

# Slide master "today" field.
$master = $p.SlideMaster
Update-DateShapeText $master.Shapes

# Every slide layout's "today" field.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShapeText $layouts.Item($li).Shapes
}

# --- Slide content edits -----------------------------------------------
$s = $p.Slides.Item(1)

# deletePerson(p) -> deleteRestaurant(r)
$deleteShape = $s.Shapes.Item(16)
$deleteRange = $deleteShape.TextFrame.TextRange
$deleteRange.Characters(1, 12).Text = "deleteRestaurant"
$deleteShape.TextFrame.TextRange.Characters(17, 3).Text = "(r)"

# Reposition/resize the text box to match the wider label (height stays
# the same, only the left edge shifts left and the width grows). The
# extra digits on Width/Height nudge the value so it still lands on the
# right EMU after the COM layer's Single-precision (float32) round trip.
$deleteShape.Left = 330.0
$deleteShape.Top = 248.0671653543307
$deleteShape.Width = 120.66007995605469
$deleteShape.Height = 16.964096069335938

# saveAddressBook(AddressBook) -> saveFoodDiary(FoodDiary)
$saveShape = $s.Shapes.Item(24)
$saveRange = $saveShape.TextFrame.TextRange
$saveRange.Characters(17, 11).Text = "FoodDiary"
$saveShape.TextFrame.TextRange.Characters(1, 15).Text = "saveFoodDiary"
